$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 897.75
$ws.Range("I40").Value = 895.5
$ws.Range("J40").Value = 900
$ws.Range("K40").Value = 895.5
$ws.Range("L40").Value = 900
$ws.Range("M40").Value = -720.5
$ws.Range("N40").Value = -1250
$ws.Range("H86").Value = 6474224
$ws.Range("I86").Value = 8463293
$ws.Range("J86").Value = 9751
$ws.Range("K86").Value = 8463293
$ws.Range("L86").Value = 9751
$ws.Range("M86").Value = -8462170
$ws.Range("N86").Value = -11997
$ws.Range("H89").Value = 6474224
$ws.Range("I89").Value = 8463293
$ws.Range("J89").Value = 9751
$ws.Range("K89").Value = 42316465
$ws.Range("L89").Value = 48755
$ws.Range("M89").Value = -42310849
$ws.Range("N89").Value = -59987
$ws.Range("H98").Value = 20455488
$ws.Range("I98").Value = 12501036
$ws.Range("J98").Value = 100000000
$ws.Range("K98").Value = 12501036
$ws.Range("L98").Value = 100000000
$ws.Range("M98").Value = -12499538
$ws.Range("N98").Value = -100002996
$ws.Range("H122").Value = 20455488
$ws.Range("I122").Value = 12501036
$ws.Range("J122").Value = 100000000
$ws.Range("K122").Value = 37503108
$ws.Range("L122").Value = 300000000
$ws.Range("M122").Value = -37500658
$ws.Range("N122").Value = -300004900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 51668.668
$ws.Range("J8").Value = 51668.668
$ws.Range("L8").Value = 51668.668
$ws.Range("N8").Value = -51956.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2257.4375
$ws.Range("I99").Value = 2222.7856
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2222.7856
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -724.7856000000002
$ws.Range("N99").Value = -5496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 33515.57
$ws.Range("I2").Value = 533
$ws.Range("K2").Value = 533
$ws.Range("M2").Value = -420
$ws.Range("H16").Value = 1073.4546
$ws.Range("I16").Value = 1001.5
$ws.Range("J16").Value = 1159.8
$ws.Range("K16").Value = 1001.5
$ws.Range("L16").Value = 1159.8
$ws.Range("M16").Value = -714.5
$ws.Range("N16").Value = -1733.8
$ws.Range("H58").Value = 4041.5557
$ws.Range("I58").Value = 1450.75
$ws.Range("K58").Value = 1450.75
$ws.Range("M58").Value = -1247.75
$ws.Range("H113").Value = 1073.4546
$ws.Range("I113").Value = 1001.5
$ws.Range("J113").Value = 1159.8
$ws.Range("K113").Value = 1001.5
$ws.Range("L113").Value = 1159.8
$ws.Range("M113").Value = 1168.5
$ws.Range("N113").Value = -5499.8
$ws.Range("H136").Value = 4041.5557
$ws.Range("I136").Value = 1450.75
$ws.Range("K136").Value = 4352.25
$ws.Range("M136").Value = -1802.25
$ws.Range("H141").Value = 39964.832
$ws.Range("J141").Value = 39964.832
$ws.Range("L141").Value = 39964.832
$ws.Range("N141").Value = -50324.832

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 825.2121
$ws.Range("I5").Value = 296.35
$ws.Range("J5").Value = 1638.8462
$ws.Range("K5").Value = 889.0500000000001
$ws.Range("L5").Value = 4916.5386
$ws.Range("M5").Value = -777.0500000000001
$ws.Range("N5").Value = -5140.5386
$ws.Range("H68").Value = 809.38464
$ws.Range("I68").Value = 454.86365
$ws.Range("J68").Value = 1069.3667
$ws.Range("K68").Value = 1364.59095
$ws.Range("L68").Value = 3208.1001
$ws.Range("M68").Value = -553.59095
$ws.Range("N68").Value = -4830.1001
$ws.Range("H71").Value = 809.38464
$ws.Range("I71").Value = 454.86365
$ws.Range("J71").Value = 1069.3667
$ws.Range("K71").Value = 4093.77285
$ws.Range("L71").Value = 9624.300300000001
$ws.Range("M71").Value = -37.77285000000029
$ws.Range("N71").Value = -17736.3003
$ws.Range("H107").Value = 340.36066
$ws.Range("I107").Value = 207.64815
$ws.Range("J107").Value = 1364.1428
$ws.Range("K107").Value = 622.94445
$ws.Range("L107").Value = 4092.4284
$ws.Range("M107").Value = 1297.05555
$ws.Range("N107").Value = -7932.428400000001
$ws.Range("H135").Value = 825.2121
$ws.Range("I135").Value = 296.35
$ws.Range("J135").Value = 1638.8462
$ws.Range("K135").Value = 2667.15
$ws.Range("L135").Value = 14749.6158
$ws.Range("M135").Value = -132.1500000000001
$ws.Range("N135").Value = -19819.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 69003
$ws.Range("J7").Value = 69003
$ws.Range("L7").Value = 69003
$ws.Range("N7").Value = -69227
$ws.Range("H8").Value = 69003
$ws.Range("J8").Value = 69003
$ws.Range("L8").Value = 69003
$ws.Range("N8").Value = -69281
$ws.Range("H97").Value = 592.8333
$ws.Range("I97").Value = 396.75
$ws.Range("J97").Value = 985
$ws.Range("K97").Value = 396.75
$ws.Range("L97").Value = 985
$ws.Range("M97").Value = 99.25
$ws.Range("N97").Value = -1977

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1472
$ws.Range("J22").Value = 1420.6
$ws.Range("L22").Value = 1420.6
$ws.Range("N22").Value = -2010.6
$ws.Range("H27").Value = 1472
$ws.Range("J27").Value = 1420.6
$ws.Range("L27").Value = 1420.6
$ws.Range("N27").Value = -1634.6
$ws.Range("H46").Value = 949.9231
$ws.Range("I46").Value = 972.1111
$ws.Range("J46").Value = 900
$ws.Range("K46").Value = 972.1111
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = -784.1111
$ws.Range("N46").Value = -1276
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H139").Value = 48516.46
$ws.Range("J139").Value = 48516.46
$ws.Range("L139").Value = 48516.46
$ws.Range("N139").Value = -58796.46

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 8913.799999999999
$ws.Range("J45").Value = 8913.799999999999
$ws.Range("L45").Value = 8913.799999999999
$ws.Range("N45").Value = -9895.799999999999
$ws.Range("H81").Value = 3060
$ws.Range("I81").Value = 1900
$ws.Range("J81").Value = 3188.889
$ws.Range("K81").Value = 3800
$ws.Range("L81").Value = 6377.778
$ws.Range("M81").Value = -2739
$ws.Range("N81").Value = -8499.778
$ws.Range("H84").Value = 3060
$ws.Range("I84").Value = 1900
$ws.Range("J84").Value = 3188.889
$ws.Range("K84").Value = 19000
$ws.Range("L84").Value = 31888.89
$ws.Range("M84").Value = -13696
$ws.Range("N84").Value = -42496.89
$ws.Range("H113").Value = 246.67857
$ws.Range("I113").Value = 241
$ws.Range("J113").Value = 263.7143
$ws.Range("K113").Value = 723
$ws.Range("L113").Value = 791.1428999999999
$ws.Range("M113").Value = 1447
$ws.Range("N113").Value = -5131.1429
